$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pspn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2111796666666667
$ws.Range("H2").Value = 0.633539
$ws.Range("I2").Value = 0.1733132136419605
$ws.Range("J2").Value = 0.1733132136419605
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03998533333333333
$ws.Range("N2").Value = 0.119956
$ws.Range("O2").Value = 0.001814551768531471
$ws.Range("P2").Value = 0.00181455176853147
$ws.Range("Q2").Value = 0.008444089364888888
$ws.Range("R2").Value = 0.07599680428399999
$ws.Range("S2").Value = 0.0003144857983238921
$ws.Range("T2").Value = 0.0003144857983238919

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pspn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2111796666666667
$ws.Range("H3").Value = 0.633539
$ws.Range("I3").Value = 0.1733132136419605
$ws.Range("J3").Value = 0.1733132136419605
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.933008
$ws.Range("N3").Value = 47.799024
$ws.Range("O3").Value = 0.7230468132755195
$ws.Range("P3").Value = 0.7230468132755195
$ws.Range("Q3").Value = 3.364727318437333
$ws.Range("R3").Value = 30.282545865936
$ws.Range("S3").Value = 0.1253135668223588
$ws.Range("T3").Value = 0.1253135668223588

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pspn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2111796666666667
$ws.Range("H4").Value = 0.633539
$ws.Range("I4").Value = 0.1733132136419605
$ws.Range("J4").Value = 0.1733132136419605
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.062935333333333
$ws.Range("N4").Value = 18.188806
$ws.Range("O4").Value = 0.275138634955949
$ws.Range("P4").Value = 0.275138634955949
$ws.Range("Q4").Value = 1.280368662714889
$ws.Range("R4").Value = 11.523317964434
$ws.Range("S4").Value = 0.04768516102127777
$ws.Range("T4").Value = 0.04768516102127777

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pspn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.007306
$ws.Range("H5").Value = 3.021918
$ws.Range("I5").Value = 0.8266867863580396
$ws.Range("J5").Value = 0.8266867863580396
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03998533333333333
$ws.Range("N5").Value = 0.119956
$ws.Range("O5").Value = 0.001814551768531471
$ws.Range("P5").Value = 0.00181455176853147
$ws.Range("Q5").Value = 0.04027746617866666
$ws.Range("R5").Value = 0.362497195608
$ws.Range("S5").Value = 0.001500065970207579
$ws.Range("T5").Value = 0.001500065970207579

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pspn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.007306
$ws.Range("H6").Value = 3.021918
$ws.Range("I6").Value = 0.8266867863580396
$ws.Range("J6").Value = 0.8266867863580396
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.933008
$ws.Range("N6").Value = 47.799024
$ws.Range("O6").Value = 0.7230468132755195
$ws.Range("P6").Value = 0.7230468132755195
$ws.Range("Q6").Value = 16.049414556448
$ws.Range("R6").Value = 144.444731008032
$ws.Range("S6").Value = 0.5977332464531607
$ws.Range("T6").Value = 0.5977332464531607

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pspn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.007306
$ws.Range("H7").Value = 3.021918
$ws.Range("I7").Value = 0.8266867863580396
$ws.Range("J7").Value = 0.8266867863580396
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.062935333333333
$ws.Range("N7").Value = 18.188806
$ws.Range("O7").Value = 0.275138634955949
$ws.Range("P7").Value = 0.275138634955949
$ws.Range("Q7").Value = 6.107231138878666
$ws.Range("R7").Value = 54.96508024990801
$ws.Range("S7").Value = 0.2274534739346713
$ws.Range("T7").Value = 0.2274534739346713
